$d = $word.ActiveDocument

# Locate the paragraph that ends with "... as good practice." - the new
# list item belongs directly after it (and before the existing, final
# empty list paragraph).
$anchorIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*as good practice.*") {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -eq $null) {
    throw "Could not find anchor paragraph containing 'as good practice.'"
}

# Insert right before the paragraph that currently follows the anchor,
# so the new paragraph lands between the anchor and it.
$nextPara = $d.Paragraphs.Item($anchorIndex + 1)
$insertRange = $nextPara.Range
$insertRange.Collapse(1)

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr>
<w:ilvl w:val="0"/>
<w:numId w:val="1"/>
</w:numPr>
</w:pPr>
<w:r><w:t xml:space="preserve">Right-click on the </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>Biggest</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t xml:space="preserve"> method and select Create Unit Tests. This creates a second project to the solution. No need to change anything.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$insertRange.InsertXML($xml)
